# Rate Verification update: refresh the FedEx ShipmentTracking numbers
# (column P, rows 2-26) on Sheet1 with a new batch of tracking numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tracking numbers, in row order (row 2 .. row 26).
$trackingNumbers = @(
  "320018714339",
  "320018714340",
  "320018714372",
  "320018714394",
  "320018714431",
  "320018714453",
  "320018714486",
  "320018714501",
  "320018714534",
  "320018714556",
  "320018714590",
  "320018714615",
  "320018714648",
  "320018714660",
  "320018714692",
  "320018714718",
  "320018714751",
  "320018714773",
  "320018714800",
  "320018714821",
  "320018714854",
  "320018714865",
  "320018714876",
  "320018714887",
  "320018714898"
)

# These values are long digit strings that Excel would otherwise interpret
# as numbers; format the target cells as Text first so they stay strings
# (matching ShipmentTracking numbers stored elsewhere in the sheet).
$startRow = 2
$endRow = $startRow + $trackingNumbers.Length - 1
$ws.Range("P${startRow}:P${endRow}").NumberFormat = "@"

for ($i = 0; $i -lt $trackingNumbers.Length; $i++) {
  $row = $startRow + $i
  $ws.Cells.Item($row, 16).Value = $trackingNumbers[$i]
}
